$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, pushing existing rows 11..38 down to 12..39.
$ws.Rows(11).Insert()

# Populate the newly inserted row 11 with the new weekly record.
$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11, 3).Value = "Ñuble"
$ws.Cells.Item(11, 4).Value = 44607
$ws.Cells.Item(11, 5).Value = 16
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100103
$ws.Cells.Item(11, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(11, 9).Value = 100103002
$ws.Cells.Item(11, 10).Value = "Ciruela"
$ws.Cells.Item(11, 11).Value = "Black Amber"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 60
$ws.Cells.Item(11, 14).Value = 11000
$ws.Cells.Item(11, 15).Value = 12000
$ws.Cells.Item(11, 16).Value = 11500
$ws.Cells.Item(11, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 639
$ws.Cells.Item(11, 20).Value = 18

# Make sure the date cell keeps the same date number format used by the other rows.
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
